$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("material_and_methods")

# Order matters: new shared strings are appended to the shared-string table
# in the order they are first written, so we match the target ordering
# (Plum, measured, EA, russian corer, roots and rhizomes included,
#  sediment not sieved, acid fumigation, gamma, organic carbon).
$ws.Range("AA4").Value = "Plum"
$ws.Range("R4").Value = "measured"
$ws.Range("U4").Value = "EA"
$ws.Range("C4").Value = "russian corer"
$ws.Range("D4").Value = "roots and rhizomes included"
$ws.Range("E4").Value = "sediment not sieved"
$ws.Range("T4").Value = "acid fumigation"
$ws.Range("X4").Value = "gamma"
$ws.Range("V4").Value = "organic carbon"

$ws.Range("Y4").Value = "gamma"
$ws.Range("L4").Value = "time approximate"
$ws.Range("H4").Value = 60
$ws.Range("I4").Value = 48
$ws.Range("S4").Value = $true

# S4/T4 pick up the formatting used elsewhere in the "data hints" rows
# (no fill applied) instead of the blank-template formatting.
$ws.Range("S1").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("T1").Copy()
$ws.Range("T4").PasteSpecial(-4122)

# Match the final view state: scrolled/zoomed out a bit, ending with AC4
# selected.
$ws.Range("Y1").Select() | Out-Null
$excel.ActiveWindow.Zoom = 109
$ws.Range("AC4").Select() | Out-Null
